$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3096085409252669
$ws1.Range("C2").Value = 0.0673076923076923
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1261261261261261
$ws1.Range("F2").Value = 0.2651515151515151
$ws1.Range("G2").Value = 0.6523297491039427
$ws1.Range("H2").Value = 0.807784911717496
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 388
$ws1.Range("K2").Value = 146
$ws1.Range("L2").Value = 0

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2734082397003745
$ws2.Range("D2").Value = 0.4294117647058823

$ws2.Range("B3").Value = 0.0673076923076923
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1261261261261261

$ws2.Range("B4").Value = 0.3096085409252669
$ws2.Range("C4").Value = 0.3096085409252669
$ws2.Range("D4").Value = 0.3096085409252669
$ws2.Range("E4").Value = 0.3096085409252669

$ws2.Range("B5").Value = 0.5336538461538461
$ws2.Range("C5").Value = 0.6367041198501873
$ws2.Range("D5").Value = 0.2777689454160042

$ws2.Range("B6").Value = 0.953531344100739
$ws2.Range("C6").Value = 0.3096085409252669
$ws2.Range("D6").Value = 0.4143014481930119

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 146
$ws3.Range("C2").Value = 388

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
